$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 105 (shifts old rows 105-170 down to 106-171)
$ws.Rows.Item(105).Insert()

# Populate the new row 105 with this week's new record
$ws.Cells.Item(105, 1).Value = 10
$ws.Cells.Item(105, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(105, 3).Value = "La Araucanía"
$ws.Cells.Item(105, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(105, 5).Value = 9
$ws.Cells.Item(105, 6).Value = 100112017
$ws.Cells.Item(105, 7).Value = "Apio"
$ws.Cells.Item(105, 8).Value = "Americana (o)"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 200
$ws.Cells.Item(105, 11).Value = 9000
$ws.Cells.Item(105, 12).Value = 9000
$ws.Cells.Item(105, 13).Value = 9000
$ws.Cells.Item(105, 14).Value = "`$/docena de matas"
$ws.Cells.Item(105, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(105, 16).Value = 1500
$ws.Cells.Item(105, 17).Value = 6
$ws.Cells.Item(105, 18).Value = "Hortaliza"
